# "Code clean up done for unwanted code in api, command and query."
#
# OOXML diff for this workbook shows:
#   1. The sheet named "Employees" is renamed to "Jul 2022".
#   2. Cosmetic/session-only bits in xl/workbook.xml (the xr:revisionPtr
#      GUID and the workbookView window rectangle) are regenerated by
#      Excel on every save and carry no document content.
#   3. Two cellXfs entries in xl/styles.xml pick up a no-op
#      applyAlignment="1" flag with no actual alignment override - a
#      save artifact that does not change how any cell looks or behaves.
#
# The only substantive, reproducible edit is the sheet rename, so that is
# what this script performs. We also nudge the window geometry via the
# Window object for completeness/parity with the author's session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename "Employees" -> "Jul 2022"
$ws.Name = "Jul 2022"

# Match the author's window placement/size recorded in the workbook view.
$win = $excel.ActiveWindow
$win.Left = 3000
$win.Top = 3000
$win.Width = 17280
$win.Height = 8880
